$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Entsoe actual-consumption series: shift existing rows 2-40 to the new
# reporting date and append rows 41-44 (PC SunEnergy additions to the forecast window)
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44)
$avals = @(5266, 5262, 5182, 5202, 5128, 5089, 5084, 5047, 5012, 5084, 5034, 5067, 5064, 5111, 5132, 5154, 5221, 5194, 5209, 5230, 5342, 5492, 5555, 5686, 5902, 6045, 6032, 6068, 6172, 6214, 6213, 6138, 6175, 6167, 6131, 5951, 5834, 5697, 5630, 5509, 5303, 5208, 5232)
$bvals = @(45799, 45799.01041666666, 45799.02083333334, 45799.03125, 45799.04166666666, 45799.05208333334, 45799.0625, 45799.07291666666, 45799.08333333334, 45799.09375, 45799.10416666666, 45799.11458333334, 45799.125, 45799.13541666666, 45799.14583333334, 45799.15625, 45799.16666666666, 45799.17708333334, 45799.1875, 45799.19791666666, 45799.20833333334, 45799.21875, 45799.22916666666, 45799.23958333334, 45799.25, 45799.26041666666, 45799.27083333334, 45799.28125, 45799.29166666666, 45799.30208333334, 45799.3125, 45799.32291666666, 45799.33333333334, 45799.34375, 45799.35416666666, 45799.36458333334, 45799.375, 45799.38541666666, 45799.39583333334, 45799.40625, 45799.41666666666, 45799.42708333334, 45799.4375)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $avals[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    if ($r -gt 40) {
        $ws.Cells.Item($r, 2).NumberFormat = $ws.Cells.Item(40, 2).NumberFormat
    }
}
